$d = $word.ActiveDocument

# Replace a whole paragraph's contents with a precise run-tree, via InsertXML
# (Find/Replace in this runtime silently merges away the leading empty <w:r/>
# run that several paragraphs carry, so we author the exact paragraph XML
# instead - this keeps that vestigial empty run intact, matching the
# canonical OOXML produced by Word's own editor.)
function Set-ParaXML($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# Paragraph 1: Heading1 title
Set-ParaXML 1 '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Million Dracula Free - Review of the Horror-Themed Slot Game</w:t></w:r>'

# Paragraphs 31-34: "What we like" bullet list
Set-ParaXML 31 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Immersive atmosphere with sound effects</w:t></w:r>'
Set-ParaXML 32 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Increasing grid size for more winning opportunities</w:t></w:r>'
Set-ParaXML 33 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Free spins feature with scatter symbols</w:t></w:r>'
Set-ParaXML 34 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Variety of horror-themed symbols</w:t></w:r>'

# Paragraph 37: "What we don't like" bullet list, second item
Set-ParaXML 37 '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Lack of music in the game</w:t></w:r>'

# Paragraphs 38-39 (bold title repeated near the end + italic closing
# summary) are the last two paragraphs in the body, immediately followed
# by the sectPr. InsertXML-ing just paragraph 39 alone (the very last
# paragraph in the story) leaves a stray trailing empty <w:p/> behind, so
# both paragraphs are replaced together in one call spanning their
# combined range - that keeps the paragraph-mark count balanced.
$p38 = $d.Paragraphs.Item(38)
$p39 = $d.Paragraphs.Item(39)
$rngTail = $d.Range($p38.Range.Start, $p39.Range.End)
$xmlTail = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
           '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Million Dracula Free - Review of the Horror-Themed Slot Game</w:t></w:r></w:p>' + `
           '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the gameplay and atmosphere of Million Dracula slot game. Play for free and win big.</w:t></w:r></w:p>' + `
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngTail.InsertXML($xmlTail)
